$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    [double]"0.9999824846980448",
    [double]"0.9990072066285755",
    [double]"0.9999999999999913",
    [double]"0.9999766506319492",
    [double]"0.9999899671105801",
    [double]"1.634976554226276e-05",
    [double]"0.0009267290336323482",
    [double]"8.853005819626963e-15",
    [double]"2.061499953517614e-05",
    [double]"1.030749977201457e-05",
    [double]"0.0002547310818618769",
    [double]"0.004043484331892824",
    [double]"1.000016814689877",
    [double]"0.004215623831067973",
    [double]"120.0425940012932",
    [double]"179.767509419835"
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
